$wb = $excel.ActiveWorkbook

# Rename the commodity label "Copper ores and concentrates" -> "Copper"
# (disaggregation of commodity Copper) on every year worksheet (cell C4).
foreach ($ws in $wb.Worksheets) {
    $label = $ws.Range("C4").Value2
    if ($label -eq "Copper ores and concentrates") {
        $ws.Range("C4").Value2 = "Copper"
    }
}

# Tiny (last-digit / ulp-level) recomputation updates to the Copper row (D4)
# that resulted from the disaggregation, for the affected year worksheets.
$wb.Worksheets.Item("2021").Range("D4").Value2 = 67537.75683225013
$wb.Worksheets.Item("2028").Range("D4").Value2 = 186790.5523943972
$wb.Worksheets.Item("2031").Range("D4").Value2 = 251326.2248433977
$wb.Worksheets.Item("2041").Range("D4").Value2 = 1011745.887999737
$wb.Worksheets.Item("2048").Range("D4").Value2 = 3900069.184158057
$wb.Worksheets.Item("2054").Range("D4").Value2 = 4192831.202238687
$wb.Worksheets.Item("2056").Range("D4").Value2 = 3479594.426315441
$wb.Worksheets.Item("2072").Range("D4").Value2 = 3907495.461863714
$wb.Worksheets.Item("2077").Range("D4").Value2 = 3528942.520219186
$wb.Worksheets.Item("2081").Range("D4").Value2 = 3161665.033672118
$wb.Worksheets.Item("2090").Range("D4").Value2 = 4104757.333033283
$wb.Worksheets.Item("2092").Range("D4").Value2 = 4170918.044394513
